# Insert two new weekly price rows for "Crimpson Seedless" grapes
# (Terminal Hortofrutícola Agro Chillán) above the existing data block,
# pushing the previously-recorded rows down by two positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 150-151; everything that was row 150.. shifts to 152..
$ws.Rows("150:151").Insert()

# --- Row 150: Crimpson Seedless / Primera ---
$ws.Range("A150").Value = 7
$ws.Range("B150").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C150").Value = "Ñuble"
$ws.Range("D150").Value = 45072
$ws.Range("E150").Value = 16
$ws.Range("F150").Value = "Fruta"
$ws.Range("G150").Value = 100109
$ws.Range("H150").Value = "Uva"
$ws.Range("I150").Value = 100109001
$ws.Range("J150").Value = "Uva"
$ws.Range("K150").Value = "Crimpson Seedless"
$ws.Range("L150").Value = "Primera"
$ws.Range("M150").Value = 50
$ws.Range("N150").Value = 12000
$ws.Range("O150").Value = 12000
$ws.Range("P150").Value = 12000
$ws.Range("Q150").Value = "`$/bandeja 18 kilos"
$ws.Range("R150").Value = "Región de O'Higgins"
$ws.Range("S150").Value = 667
$ws.Range("T150").Value = 18

# --- Row 151: Crimpson Seedless / Segunda ---
$ws.Range("A151").Value = 7
$ws.Range("B151").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C151").Value = "Ñuble"
$ws.Range("D151").Value = 45072
$ws.Range("E151").Value = 16
$ws.Range("F151").Value = "Fruta"
$ws.Range("G151").Value = 100109
$ws.Range("H151").Value = "Uva"
$ws.Range("I151").Value = 100109001
$ws.Range("J151").Value = "Uva"
$ws.Range("K151").Value = "Crimpson Seedless"
$ws.Range("L151").Value = "Segunda"
$ws.Range("M151").Value = 60
$ws.Range("N151").Value = 10000
$ws.Range("O151").Value = 10000
$ws.Range("P151").Value = 10000
$ws.Range("Q151").Value = "`$/bandeja 18 kilos"
$ws.Range("R151").Value = "Región de O'Higgins"
$ws.Range("S151").Value = 556
$ws.Range("T151").Value = 18
